$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 22 ("serviceName"), shifting the
# remaining "Service" table rows down by one.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row with the new "contactMail" field.
$ws.Range("A22").Value = "contactMail"
$ws.Range("B22").Value = "string"
